$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 62.41592166666667
$ws.Range("H2").Value = 187.247765
$ws.Range("I2").Value = 0.1654944774607044
$ws.Range("J2").Value = 0.1654944774607044
$ws.Range("M2").Value = 8.820647333333334
$ws.Range("N2").Value = 26.461942
$ws.Range("O2").Value = 0.06415146660411865
$ws.Range("P2").Value = 0.06415146660411865
$ws.Range("Q2").Value = 550.5488330066256
$ws.Range("R2").Value = 4954.93949705963
$ws.Range("S2").Value = 0.01061671344398645
$ws.Range("T2").Value = 0.01061671344398645
$ws.Range("G3").Value = 62.41592166666667
$ws.Range("H3").Value = 187.247765
$ws.Range("I3").Value = 0.1654944774607044
$ws.Range("J3").Value = 0.1654944774607044
$ws.Range("O3").Value = 0.3979101621202897
$ws.Range("P3").Value = 0.3979101621202898
$ws.Range("Q3").Value = 3414.870882823092
$ws.Range("R3").Value = 43316.78728617237
$ws.Range("S3").Value = 0.06585193435640155
$ws.Range("T3").Value = 0.06585193435640155
$ws.Range("G4").Value = 62.41592166666667
$ws.Range("H4").Value = 187.247765
$ws.Range("I4").Value = 0.1654944774607044
$ws.Range("J4").Value = 0.1654944774607044
$ws.Range("M4").Value = 21.90816333333333
$ws.Range("N4").Value = 65.72449
$ws.Range("O4").Value = 0.1593353362087987
$ws.Range("P4").Value = 0.1593353362087987
$ws.Range("Q4").Value = 1367.418206473872
$ws.Range("R4").Value = 12306.76385826485
$ws.Range("S4").Value = 0.0263691182069008
$ws.Range("T4").Value = 0.0263691182069008
$ws.Range("G5").Value = 62.41592166666667
$ws.Range("H5").Value = 187.247765
$ws.Range("I5").Value = 0.1654944774607044
$ws.Range("J5").Value = 0.1654944774607044
$ws.Range("M5").Value = 52.056859
$ws.Range("N5").Value = 156.170577
$ws.Range("O5").Value = 0.3786030350667928
$ws.Range("P5").Value = 0.3786030350667929
$ws.Range("Q5").Value = 3249.176833556712
$ws.Range("R5").Value = 29242.5915020104
$ws.Range("S5").Value = 0.06265671145341564
$ws.Range("T5").Value = 0.06265671145341566
$ws.Range("I6").Value = 0.4369365253446571
$ws.Range("J6").Value = 0.436936525344657
$ws.Range("M6").Value = 8.820647333333334
$ws.Range("N6").Value = 26.461942
$ws.Range("O6").Value = 0.06415146660411865
$ws.Range("P6").Value = 0.06415146660411865
$ws.Range("Q6").Value = 1453.552395327445
$ws.Range("R6").Value = 13081.97155794701
$ws.Range("S6").Value = 0.02803011891376741
$ws.Range("T6").Value = 0.02803011891376741
$ws.Range("I7").Value = 0.4369365253446571
$ws.Range("J7").Value = 0.436936525344657
$ws.Range("O7").Value = 0.3979101621202897
$ws.Range("P7").Value = 0.3979101621202898
$ws.Range("S7").Value = 0.1738614836361686
$ws.Range("T7").Value = 0.1738614836361686
$ws.Range("I8").Value = 0.4369365253446571
$ws.Range("J8").Value = 0.436936525344657
$ws.Range("M8").Value = 21.90816333333333
$ws.Range("N8").Value = 65.72449
$ws.Range("O8").Value = 0.1593353362087987
$ws.Range("P8").Value = 0.1593353362087987
$ws.Range("Q8").Value = 3610.241072676175
$ws.Range("R8").Value = 32492.16965408557
$ws.Range("S8").Value = 0.06961942816769522
$ws.Range("T8").Value = 0.06961942816769522
$ws.Range("I9").Value = 0.4369365253446571
$ws.Range("J9").Value = 0.436936525344657
$ws.Range("M9").Value = 52.056859
$ws.Range("N9").Value = 156.170577
$ws.Range("O9").Value = 0.3786030350667928
$ws.Range("P9").Value = 0.3786030350667929
$ws.Range("Q9").Value = 8578.437526543563
$ws.Range("R9").Value = 77205.93773889205
$ws.Range("S9").Value = 0.1654254946270258
$ws.Range("T9").Value = 0.1654254946270258
$ws.Range("G10").Value = 57.486235
$ws.Range("H10").Value = 172.458705
$ws.Range("I10").Value = 0.1524235190071549
$ws.Range("J10").Value = 0.1524235190071549
$ws.Range("M10").Value = 8.820647333333334
$ws.Range("N10").Value = 26.461942
$ws.Range("O10").Value = 0.06415146660411865
$ws.Range("P10").Value = 0.06415146660411865
$ws.Range("Q10").Value = 507.0658054561234
$ws.Range("R10").Value = 4563.592249105111
$ws.Range("S10").Value = 0.009778192289269742
$ws.Range("T10").Value = 0.009778192289269742
$ws.Range("G11").Value = 57.486235
$ws.Range("H11").Value = 172.458705
$ws.Range("I11").Value = 0.1524235190071549
$ws.Range("J11").Value = 0.1524235190071549
$ws.Range("O11").Value = 0.3979101621202897
$ws.Range("P11").Value = 0.3979101621202898
$ws.Range("Q11").Value = 3145.160158220725
$ws.Range("R11").Value = 28306.44142398653
$ws.Range("S11").Value = 0.06065086715908207
$ws.Range("T11").Value = 0.06065086715908208
$ws.Range("G12").Value = 57.486235
$ws.Range("H12").Value = 172.458705
$ws.Range("I12").Value = 0.1524235190071549
$ws.Range("J12").Value = 0.1524235190071549
$ws.Range("M12").Value = 21.90816333333333
$ws.Range("N12").Value = 65.72449
$ws.Range("O12").Value = 0.1593353362087987
$ws.Range("P12").Value = 0.1593353362087987
$ws.Range("Q12").Value = 1259.417825798383
$ws.Range("R12").Value = 11334.76043218545
$ws.Range("S12").Value = 0.02428645264713324
$ws.Range("T12").Value = 0.02428645264713324
$ws.Range("G13").Value = 57.486235
$ws.Range("H13").Value = 172.458705
$ws.Range("I13").Value = 0.1524235190071549
$ws.Range("J13").Value = 0.1524235190071549
$ws.Range("M13").Value = 52.056859
$ws.Range("N13").Value = 156.170577
$ws.Range("O13").Value = 0.3786030350667928
$ws.Range("P13").Value = 0.3786030350667929
$ws.Range("Q13").Value = 2992.552829835865
$ws.Range("R13").Value = 26932.97546852278
$ws.Range("S13").Value = 0.05770800691166983
$ws.Range("T13").Value = 0.05770800691166984
$ws.Range("G14").Value = 92.45614233333333
$ws.Range("H14").Value = 277.368427
$ws.Range("I14").Value = 0.2451454781874835
$ws.Range("J14").Value = 0.2451454781874835
$ws.Range("M14").Value = 8.820647333333334
$ws.Range("N14").Value = 26.461942
$ws.Range("O14").Value = 0.06415146660411865
$ws.Range("P14").Value = 0.06415146660411865
$ws.Range("Q14").Value = 815.5230253228037
$ws.Range("R14").Value = 7339.707227905235
$ws.Range("S14").Value = 0.01572644195709505
$ws.Range("T14").Value = 0.01572644195709505
$ws.Range("G15").Value = 92.45614233333333
$ws.Range("H15").Value = 277.368427
$ws.Range("I15").Value = 0.2451454781874835
$ws.Range("J15").Value = 0.2451454781874835
$ws.Range("O15").Value = 0.3979101621202897
$ws.Range("P15").Value = 0.3979101621202898
$ws.Range("Q15").Value = 5058.417467235148
$ws.Range("R15").Value = 45525.75720511633
$ws.Range("S15").Value = 0.09754587696863752
$ws.Range("T15").Value = 0.09754587696863753
$ws.Range("G16").Value = 92.45614233333333
$ws.Range("H16").Value = 277.368427
$ws.Range("I16").Value = 0.2451454781874835
$ws.Range("J16").Value = 0.2451454781874835
$ws.Range("M16").Value = 21.90816333333333
$ws.Range("N16").Value = 65.72449
$ws.Range("O16").Value = 0.1593353362087987
$ws.Range("P16").Value = 0.1593353362087987
$ws.Range("Q16").Value = 2025.544267408581
$ws.Range("R16").Value = 18229.89840667723
$ws.Range("S16").Value = 0.03906033718706941
$ws.Range("T16").Value = 0.03906033718706941
$ws.Range("G17").Value = 92.45614233333333
$ws.Range("H17").Value = 277.368427
$ws.Range("I17").Value = 0.2451454781874835
$ws.Range("J17").Value = 0.2451454781874835
$ws.Range("M17").Value = 52.056859
$ws.Range("N17").Value = 156.170577
$ws.Range("O17").Value = 0.3786030350667928
$ws.Range("P17").Value = 0.3786030350667929
$ws.Range("Q17").Value = 4812.976365130264
$ws.Range("R17").Value = 43316.78728617237
$ws.Range("S17").Value = 0.09281282207468153
$ws.Range("T17").Value = 0.09281282207468154
